# Apply print-request-logsheet updates:
#  - FOCUS a4c 50 scale row (21) gets a completed date of 30-05-2018 in column B
#  - FOCUS LAX 50 scale completion count corrected from 12 to 18
#  - FOCUS SAX 50 scale completion count corrected from 11 to 15
#  - New request added for Tracheal Rings (row 34)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the date completed for the "FOCUS a4c 50 scale" request
$ws.Range("B21").Value = "30-05-2018"

# Update completion counts
$ws.Range("I23").Value = "18 complete"
$ws.Range("I25").Value = "15 complete"

# New row for the Tracheal Rings print request
$ws.Range("A34").Value = "31-05-2018"
$ws.Range("B34").Value = "31-05-2018"
$ws.Range("C34").Value = "Tracheal Rings"
$ws.Range("D34").Value = 15
$ws.Range("E34").Value = "PLA"
$ws.Range("F34").Value = 2
$ws.Range("G34").Value = 20
$ws.Range("H34").Value = 0.2
$ws.Range("I34").Value = "NA"

# Match row height that Excel derives for row 21 once it holds a value in column B
$ws.Rows.Item(21).RowHeight = 13.8

# Keep the view in a sensible state around the newly-edited rows
$ws.Range("I23").Select()
